$d = $word.ActiveDocument

# Locate the start paragraph ("Problem space...") and the end paragraph (last
# bullet, "would they be willing to adopt...") by searching paragraph text,
# since relying on fixed paragraph indices would be fragile.
$startIndex = -1
$endIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIndex -eq -1 -and $t.Contains("Problem space")) {
        $startIndex = $i
    }
    if ($t.Contains("would they be willing to adopt a new device specific to managing the condition?")) {
        $endIndex = $i
    }
}

if ($startIndex -eq -1 -or $endIndex -eq -1) {
    throw "Could not locate target paragraphs (start=$startIndex end=$endIndex)"
}

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Problem space</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">Considering the increasing number of smartphones users, we are interested in the management of health conditions using digital devices. In other words, we would like to understand how specialized digital devices and/or apps can help people manage their health conditions. </w:t></w:r></w:p><w:p><w:r><w:t>During the process of exploring the exact problem, we hope to find out what are the issues users have when they are monitoring / managing diabetes.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>User Group:</w:t></w:r><w:r><w:t xml:space="preserve"> the user group that we are interested in is people who are suffering from diabetes (both types 1 and 2). For many people, this disease requires constant monitoring of health factors such as blood glucose levels.</w:t></w:r></w:p><w:p><w:r><w:t>Monitoring is required to keep symptoms under control and reduce the risk of progression and complications. Many sufferers report ‘burn-out’ and stress from managing their condition.</w:t></w:r></w:p><w:p><w:r><w:t>Seniors in particular are much less likely to embrace technology to facilitate the management of their health condition</w:t></w:r><w:r><w:t xml:space="preserve">, so we would like to understand their difficulties with technology and create a helpful interface to help them. </w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">High Level </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Research plan:</w:t></w:r><w:r><w:t xml:space="preserve"> We intend to seek out research participants from local (Toronto) support groups. For instance: the Canadian Diabetes Association (CDA) organizes local support groups. Another possible source of users may be the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Banting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &amp; Best Diabetes Center (BBDC) at the University of Toronto.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The purpose of our research is to obtain a specific understanding of how seniors manage their diabetes (i.e. tools and techniques), what further steps they need to take to manage their condition, and what prevents them from doing so.</w:t></w:r></w:p><w:p><w:r><w:t>The primary research method will be through interviews, as conversations are likely to provide greater insight and allow the asking of open questions and following up on the answers. It is possible that a survey may be developed from the interview results in order to confirm information from a broader sampling.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p><w:p><w:r><w:t>Interview questions will explore the following topics:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>How</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> do they currently manage their condition</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>What</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> issues do they face with managing their condition</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>What</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> do they need to manage</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>What</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> tools do they currently have</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>How</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> much do they use them</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>What</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> are their (tools’) strengths and shortcomings</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>In</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> general, what technology devices do they use (for any purpose)</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Would</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> they be willing to adopt a new device specific to managing the condition?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)

Write-Host "Replaced paragraphs $startIndex through $endIndex"
